# "Fix SDG data + wbcodesUN + create briefs"
#
# The sheet's AutoFilter previously filtered column C ("source") down to
# rows where source = "ILO". Re-point the filter at column B ("code") and
# select the MNCH_* indicator rows instead - this both rewrites the
# <autoFilter> definition and (as a side effect, exactly like real Excel)
# re-hides/un-hides every data row so only the matching rows stay visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:L144")

# Turn off the existing filter on column 3 (C = "source"), which was
# previously restricting the view to val="ILO".
$dataRange.AutoFilter(3)

# Turn on a values-filter on column 2 (B = "code") restricted to the
# MNCH_* indicator codes.
$mnchCodes = @( `
    "MNCH_ANC1", `
    "MNCH_ANC4", `
    "MNCH_BIRTH18", `
    "MNCH_DEMAND_FP", `
    "MNCH_DIARCARE", `
    "MNCH_ITN", `
    "MNCH_ITNPREG", `
    "MNCH_MLRACT", `
    "MNCH_MLRCARE", `
    "MNCH_MLRDIAG", `
    "MNCH_ORS", `
    "MNCH_PNCMOM", `
    "MNCH_PNCNB", `
    "MNCH_PNEUCARE", `
    "MNCH_SAB" `
)
$dataRange.AutoFilter(2, $mnchCodes, 7)

# Leave the selection where the author last was when they saved.
$ws.Range("B110").Select()
